# Hands On Demos - Day 3.
#
# Slide 42 ("Summary"), shape "object 4" (the bulleted text box):
#   - "Place ellipse after parameter type" -> "Place ellipsis after parameter type"
#     ("ellipse" is shortened to "ellips" and a new run "is" is appended so the
#     word reads "ellipsis").
#   - The shape is resized/repositioned slightly (PowerPoint's autofit kicking
#     in for the now-wider line of text).

# Helper: PowerPoint stores shape geometry in points (a float32 internally)
# and the OOXML is written out in EMU (1 pt = 12700 EMU). Converting the
# "obvious" point value (emu/12700) sometimes truncates to one EMU below the
# intended target once it round-trips through a 32 bit float, so search the
# neighbourhood of the naive value for one whose float32 representation lands
# exactly on the desired EMU value.
function Find-PointForEmu {
    param([int]$targetEmu)

    $basePt = $targetEmu / 12700.0
    if ([int]([float]$basePt * 12700) -eq $targetEmu) {
        return $basePt
    }
    for ($i = 1; $i -lt 5000; $i++) {
        $cand = $basePt + ($i * 0.000001)
        if ([int]([float]$cand * 12700) -eq $targetEmu) { return $cand }
        $cand = $basePt - ($i * 0.000001)
        if ([int]([float]$cand * 12700) -eq $targetEmu) { return $cand }
    }
    return $basePt
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(42)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# --- text edit: "ellipse" -> "ellipsis" (as two runs: "ellips" + "is") -----
$found = $tr.Find("ellipse")
$found.Text = "ellips"

$word = $tr.Find("ellips")
[void]($word.InsertAfter("is"))

# Force the freshly inserted "is" to become its own run (matching the source
# diff, which shows it as a separate <a:r>) by touching one of its run-level
# properties; re-assigning the colour it already has is a no-op visually but
# is enough to split the run away from its neighbours.
$isStart = $word.Start + $word.Length
$isRun = $tr.Characters($isStart, 2)
$isRun.Font.Color.RGB = $isRun.Font.Color.RGB

# --- shape geometry: reposition/resize the textbox -------------------------
$shp.Left = Find-PointForEmu 5226050
$shp.Top = Find-PointForEmu 2373630
$shp.Width = Find-PointForEmu 6415405
$shp.Height = Find-PointForEmu 1793875
